$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a text value to a cell while preserving a plain (unstyled)
# text cell - temporarily switch to text format so Excel does not
# auto-coerce numeric-looking strings (e.g. "1.00") into numbers, then
# restore the default "Normal" style so no stray formatting is left behind.
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

# Row 2
Set-TextValue $ws.Range('D2') '66.035.48'
Set-TextValue $ws.Range('E2') '  -1.65%  '

# Row 3
Set-TextValue $ws.Range('D3') '3.416.89'
Set-TextValue $ws.Range('E3') '  -1.21%  '

# Row 4
Set-TextValue $ws.Range('E4') '  +0.03%  '

# Row 5
Set-TextValue $ws.Range('D5') '580.38'
Set-TextValue $ws.Range('E5') '  -2.21%  '

# Row 6
Set-TextValue $ws.Range('D6') '172.69'
Set-TextValue $ws.Range('E6') '  -3.71%  '

# Row 7
Set-TextValue $ws.Range('E7') '  +0.08%  '

# Row 8
Set-TextValue $ws.Range('D8') '0.587'
Set-TextValue $ws.Range('E8') '  -3.82%  '

# Row 9
Set-TextValue $ws.Range('D9') '3.421.18'
Set-TextValue $ws.Range('E9') '  -1.07%  '

# Row 10
Set-TextValue $ws.Range('D10') '0.129'
Set-TextValue $ws.Range('E10') '  -6.71%  '

# Row 11
Set-TextValue $ws.Range('D11') '6.85'
Set-TextValue $ws.Range('E11') '  -1.54%  '

# Row 12
Set-TextValue $ws.Range('D12') '0.406'
Set-TextValue $ws.Range('E12') '  -5.40%  '

# Row 13
Set-TextValue $ws.Range('D13') '4.027.70'
Set-TextValue $ws.Range('E13') '  -0.78%  '

# Row 14
Set-TextValue $ws.Range('E14') '  -0.80%  '

# Row 15
Set-TextValue $ws.Range('D15') '29.65'
Set-TextValue $ws.Range('E15') '  -7.29%  '

# Row 16
Set-TextValue $ws.Range('D16') '66.070.66'
Set-TextValue $ws.Range('E16') '  -1.61%  '

# Row 17
Set-TextValue $ws.Range('D17') '0.0000169'
Set-TextValue $ws.Range('E17') '  -4.49%  '

# Row 18
Set-TextValue $ws.Range('D18') '3.430.34'
Set-TextValue $ws.Range('E18') '  -0.89%  '

# Row 19
Set-TextValue $ws.Range('D19') '5.84'
Set-TextValue $ws.Range('E19') '  -5.72%  '

# Row 20
Set-TextValue $ws.Range('D20') '13.58'
Set-TextValue $ws.Range('E20') '  -3.89%  '

# Row 21
Set-TextValue $ws.Range('D21') '370.66'
Set-TextValue $ws.Range('E21') '  -5.63%  '

# Row 22
Set-TextValue $ws.Range('D22') '7.65'
Set-TextValue $ws.Range('E22') '  -3.24%  '

# Row 23
Set-TextValue $ws.Range('E23') '  +0.04%  '

# Row 24
Set-TextValue $ws.Range('B24') 'LEO'
Set-TextValue $ws.Range('C24') 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
Set-TextValue $ws.Range('D24') '5.73'
Set-TextValue $ws.Range('E24') '  -0.80%  '

# Row 25
Set-TextValue $ws.Range('B25') 'Litecoin'
Set-TextValue $ws.Range('C25') 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
Set-TextValue $ws.Range('D25') '71.51'
Set-TextValue $ws.Range('E25') '  -0.11%  '

# Row 26
Set-TextValue $ws.Range('D26') '0.526'
Set-TextValue $ws.Range('E26') '  -2.23%  '

# Row 27
Set-TextValue $ws.Range('E27') '  -2.85%  '

# Row 28
Set-TextValue $ws.Range('D28') '9.54'
Set-TextValue $ws.Range('E28') '  -7.62%  '

# Row 29
Set-TextValue $ws.Range('D29') '0.176'
Set-TextValue $ws.Range('E29') '  +1.05%  '

# Row 30
Set-TextValue $ws.Range('D30') '1.00'
Set-TextValue $ws.Range('E30') '  +0.08%  '

# Row 31
Set-TextValue $ws.Range('D31') '23.76'
Set-TextValue $ws.Range('E31') '  +1.29%  '

# Row 32
Set-TextValue $ws.Range('D32') '5.73'
Set-TextValue $ws.Range('E32') '  -6.25%  '

# Row 33
Set-TextValue $ws.Range('D33') '1.97'
Set-TextValue $ws.Range('E33') '  -3.83%  '

# Row 34
Set-TextValue $ws.Range('D34') '0.999'
Set-TextValue $ws.Range('E34') '  -0.02%  '

# Row 35
Set-TextValue $ws.Range('D35') '1.28'
Set-TextValue $ws.Range('E35') '  -7.90%  '

# Row 36
Set-TextValue $ws.Range('D36') '7.00'
Set-TextValue $ws.Range('E36') '  -4.20%  '

# Row 37
Set-TextValue $ws.Range('D37') '1.52'
Set-TextValue $ws.Range('E37') '  -3.50%  '

# Row 38
Set-TextValue $ws.Range('D38') '160.71'
Set-TextValue $ws.Range('E38') '  +0.13%  '

# Row 39
Set-TextValue $ws.Range('D39') '29.04'
Set-TextValue $ws.Range('E39') '  +11.44%  '

# Row 40
Set-TextValue $ws.Range('D40') '0.882'
Set-TextValue $ws.Range('E40') '  +0.69%  '

# Row 41
Set-TextValue $ws.Range('D41') '2.63'
Set-TextValue $ws.Range('E41') '  -5.78%  '

# Row 42
Set-TextValue $ws.Range('D42') '1.73'
Set-TextValue $ws.Range('E42') '  -6.91%  '

# Row 43
Set-TextValue $ws.Range('D43') '2.692.50'
Set-TextValue $ws.Range('E43') '  -2.23%  '

# Row 44
Set-TextValue $ws.Range('D44') '4.37'
Set-TextValue $ws.Range('E44') '  -6.05%  '

# Row 45
Set-TextValue $ws.Range('D45') '6.25'
Set-TextValue $ws.Range('E45') '  -6.97%  '

# Row 46
Set-TextValue $ws.Range('D46') '0.0677'
Set-TextValue $ws.Range('E46') '  -5.84%  '

# Row 47
Set-TextValue $ws.Range('D47') '40.05'
Set-TextValue $ws.Range('E47') '  -3.07%  '

# Row 48
Set-TextValue $ws.Range('D48') '0.0286'
Set-TextValue $ws.Range('E48') '  -3.79%  '

# Row 49
Set-TextValue $ws.Range('D49') '23.79'
Set-TextValue $ws.Range('E49') '  -8.91%  '

# Row 50
Set-TextValue $ws.Range('D50') '302.75'
Set-TextValue $ws.Range('E50') '  -6.62%  '

# Row 51
Set-TextValue $ws.Range('B51') 'Stellar'
Set-TextValue $ws.Range('C51') 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
Set-TextValue $ws.Range('D51') '0.100'
Set-TextValue $ws.Range('E51') '  -4.33%  '
